$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.648.07'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.591.92'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.511'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0615'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.244'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.817.43'
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").Value = '1.600.71'
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '26.622.65'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '206.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("E23").Value = '  -2.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0503'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.664'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '1.278.75'
$ws.Range("E35").Value = '  -3.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '1.728.73'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.906'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("E48").Value = '  -0.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.26%  '
